$d = $word.ActiveDocument

# Locate the existing list item "Paralelizar segmentos de código" and add a
# new list item "Desarrollo paralelo del algoritmo" right after it (before
# "Estudiar e incorporar hostfile"), inheriting the same paragraph/run
# formatting (Monserrat font, numbered-list numId 3) via InsertParagraphAfter.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Paralelizar segmentos de c.digo") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph 'Paralelizar segmentos de código'"
}

$anchor = $d.Paragraphs.Item($anchorIndex)
$anchor.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = "Desarrollo paralelo del algoritmo"
